$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(1, 1).Value = 'running knee pads'
$ws.Cells.Item(2, 1).Value = 'youth yoga pants'
$ws.Cells.Item(3, 1).Value = 'boys lacrosse pants'
$ws.Cells.Item(4, 1).Value = 'knee length baseball pants'
$ws.Cells.Item(5, 1).Value = 'youth compression knee pad sleeve'
$ws.Cells.Item(6, 1).Value = 'yoga pad for knees'
$ws.Cells.Item(7, 1).Value = 'recovery leggings men'
$ws.Cells.Item(8, 1).Value = 'hex leg sleeves basketball'
$ws.Cells.Item(9, 1).Value = 'adult basketball'
$ws.Cells.Item(10, 1).Value = 'knee pads boys youth'
$ws.Cells.Item(11, 1).Value = 'wrestling leggings for men'
$ws.Cells.Item(12, 1).Value = 'softball sliding shorts youth'
$ws.Cells.Item(13, 1).Value = 'leggings sports'
$ws.Cells.Item(14, 1).Value = 'basketball padded shorts'
$ws.Cells.Item(15, 1).Value = 'best yoga knee pad'
$ws.Cells.Item(16, 1).Value = 'mens yoga pants tight'
$ws.Cells.Item(17, 1).Value = 'baseball tights'
$ws.Cells.Item(18, 1).Value = 'volleyball catcher'
$ws.Cells.Item(19, 1).Value = 'goalkeeper pads'
$ws.Cells.Item(20, 1).Value = 'youth baseball pants knee'
$ws.Cells.Item(21, 1).Value = 'paintball knee pads'
$ws.Cells.Item(22, 1).Value = 'youth compression tights boys'
$ws.Cells.Item(23, 1).Value = 'compression tights for boys basketball'
$ws.Cells.Item(24, 1).Value = 'breathable knee pads'
$ws.Cells.Item(25, 1).Value = 'mens spandex pants'
$ws.Cells.Item(26, 1).Value = 'soccer pads boys'
$ws.Cells.Item(27, 1).Value = 'hockey knee pads youth'
$ws.Cells.Item(28, 1).Value = 'padded shorts basketball'
$ws.Cells.Item(29, 1).Value = 'sliding shorts baseball youth'
$ws.Cells.Item(30, 1).Value = 'running pads'
$ws.Cells.Item(31, 1).Value = 'rodillera de basketball'
$ws.Cells.Item(32, 1).Value = 'leg tights for men'
$ws.Cells.Item(33, 1).Value = 'mens baseball pants black'
$ws.Cells.Item(34, 1).Value = 'sliding shorts boys'
$ws.Cells.Item(35, 1).Value = 'baseball pants knee'
$ws.Cells.Item(36, 1).Value = 'adult softball pants'
$ws.Cells.Item(37, 1).Value = 'catcher knee support'
$ws.Cells.Item(38, 1).Value = 'cycling pads for men'
$ws.Cells.Item(39, 1).Value = 'youth football girdle with knee pads'
$ws.Cells.Item(40, 1).Value = 'best knee pads for basketball'
$ws.Cells.Item(41, 1).Value = 'knee pads yoga'
$ws.Cells.Item(42, 1).Value = 'compression pads'
$ws.Cells.Item(43, 1).Value = 'knee pants boys'
$ws.Cells.Item(44, 1).Value = 'basketball chart'
$ws.Cells.Item(45, 1).Value = 'girl basketball knee pads'
$ws.Cells.Item(46, 1).Value = 'knee pads for basketball girls'
$ws.Cells.Item(47, 1).Value = 'knee pads men'
$ws.Cells.Item(48, 1).Value = 'basketball compression pants boys'
$ws.Cells.Item(49, 1).Value = 'mens compression knee'
$ws.Cells.Item(50, 1).Value = 'knee pads by design'
$ws.Cells.Item(51, 1).Value = 'good knee pads'
$ws.Cells.Item(52, 1).Value = 'wrestling knee pad youth'
$ws.Cells.Item(53, 1).Value = 'mens work pants with knee pad'
$ws.Cells.Item(54, 1).Value = 'kneepad pants'
$ws.Cells.Item(55, 1).Value = 'softball girls sliding shorts'
$ws.Cells.Item(56, 1).Value = 'men yoga capri'
$ws.Cells.Item(57, 1).Value = 'compression shorts lacrosse'
$ws.Cells.Item(58, 1).Value = 'black knee pads for basketball'
$ws.Cells.Item(59, 1).Value = 'mens 3/4 compression pants'
$ws.Cells.Item(60, 1).Value = 'softball sliding shorts youth girls'
$ws.Cells.Item(61, 1).Value = 'basketball knee pads girls'
$ws.Cells.Item(62, 1).Value = 'sliding baseball shorts'
$ws.Cells.Item(63, 1).Value = 'knee pads for big men'
$ws.Cells.Item(64, 1).Value = 'sleeve knee pads basketball'
$ws.Cells.Item(65, 1).Value = 'youth padded compression shorts basketball'
$ws.Cells.Item(66, 1).Value = 'mens compression recovery pants'
$ws.Cells.Item(67, 1).Value = 'arthritis friendly yoga'
$ws.Cells.Item(68, 1).Value = 'youth compression pants for boys'
$ws.Cells.Item(69, 1).Value = 'sports leggings for men'
$ws.Cells.Item(70, 1).Value = 'long knee pads basketball'
$ws.Cells.Item(71, 1).Value = 'best knee pads for yoga'
$ws.Cells.Item(72, 1).Value = 'rodillera basketball'
$ws.Cells.Item(73, 1).Value = 'basketball compression shorts'
$ws.Cells.Item(74, 1).Value = 'mountain biking knee pads for men'
$ws.Cells.Item(75, 1).Value = 'boys athletic tights youth'
$ws.Cells.Item(76, 1).Value = 'youth boys compression tights'
$ws.Cells.Item(77, 1).Value = 'boys compression pants football'
$ws.Cells.Item(78, 1).Value = 'padded sliding shorts mens'
$ws.Cells.Item(79, 1).Value = 'softball compression shorts'
$ws.Cells.Item(80, 1).Value = 'youth boys compression leggings'
$ws.Cells.Item(81, 1).Value = 'fitness gear floor guard'
$ws.Cells.Item(82, 1).Value = 'boys compression tights basketball'
$ws.Cells.Item(83, 1).Value = 'padded football pants youth'
$ws.Cells.Item(84, 1).Value = 'football girdle with knee pads'
$ws.Cells.Item(85, 1).Value = 'youth girls sliding shorts softball'
$ws.Cells.Item(86, 1).Value = 'construction need pads'
$ws.Cells.Item(87, 1).Value = 'softball sliding shorts for girls'
$ws.Cells.Item(88, 1).Value = 'softball pants men'
$ws.Cells.Item(89, 1).Value = 'pads for running'
$ws.Cells.Item(90, 1).Value = 'big and tall compression pants'
$ws.Cells.Item(91, 1).Value = 'hex pad knee sleeve'
$ws.Cells.Item(92, 1).Value = 'adult softball'
$ws.Cells.Item(93, 1).Value = 'compression knee sleeve men basketball'
$ws.Cells.Item(94, 1).Value = 'volleyball hand protector'
$ws.Cells.Item(95, 1).Value = 'yoga knee pad'
$ws.Cells.Item(96, 1).Value = 'knee sleeve wrestling youth'
$ws.Cells.Item(97, 1).Value = 'compression recovery pants'
$ws.Cells.Item(98, 1).Value = 'below knee leggings'
$ws.Cells.Item(99, 1).Value = 'basketball padded compression'
$ws.Cells.Item(100, 1).Value = 'knee length basketball shorts for women'
